$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 30; $r++) {
    $ws.Cells.Item($r, 15).Value = "2022-07-22 20:57:25"
}
